$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7471
$ws.Range("K3").Value = 7728
$ws.Range("K4").Value = 1626
$ws.Range("K5").Value = 550
$ws.Range("K6").Value = 8607
$ws.Range("K7").Value = 25982

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 225
$ws.Range("K4").Value = 91
$ws.Range("K8").Value = 1694
$ws.Range("K9").Value = 126
$ws.Range("K11").Value = 466
$ws.Range("K14").Value = 122
$ws.Range("K15").Value = 266
$ws.Range("K19").Value = 752
$ws.Range("K25").Value = 120
$ws.Range("K29").Value = 1432
$ws.Range("K31").Value = 305
$ws.Range("K36").Value = 331
$ws.Range("K42").Value = 960
$ws.Range("K43").Value = 212
$ws.Range("K45").Value = 37
$ws.Range("K47").Value = 178
$ws.Range("K49").Value = 146
$ws.Range("K50").Value = 120
$ws.Range("K51").Value = 333
$ws.Range("K52").Value = 673
$ws.Range("K54").Value = 513
$ws.Range("K55").Value = 286
$ws.Range("K60").Value = 153
$ws.Range("K63").Value = 73
$ws.Range("K67").Value = 1014
$ws.Range("K72").Value = 126
$ws.Range("K78").Value = 320
$ws.Range("K79").Value = 637
$ws.Range("K85").Value = 1188
$ws.Range("K86").Value = 160
$ws.Range("K88").Value = 277
$ws.Range("K89").Value = 389
$ws.Range("K91").Value = 307
$ws.Range("K94").Value = 346
$ws.Range("K95").Value = 430
$ws.Range("K97").Value = 210
$ws.Range("K98").Value = 138
$ws.Range("K99").Value = 437
$ws.Range("K101").Value = 25982

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 471
$ws.Range("K4").Value = 97
$ws.Range("K6").Value = 569
$ws.Range("K7").Value = 1694

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 146
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 430

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 185
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 437

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 92
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 286
$ws.Range("K7").Value = 1014

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 78
$ws.Range("K6").Value = 281
$ws.Range("K7").Value = 513

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K5").Value = 36
$ws.Range("K6").Value = 419
$ws.Range("K7").Value = 1432

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 24
$ws.Range("K7").Value = 752

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 79
$ws.Range("K6").Value = 177

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 284
$ws.Range("K4").Value = 44
$ws.Range("K6").Value = 357
$ws.Range("K7").Value = 960

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 95
$ws.Range("K7").Value = 320

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 286

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 144
$ws.Range("K7").Value = 307

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K4").Value = 43
$ws.Range("K6").Value = 166
$ws.Range("K7").Value = 637

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 199
$ws.Range("K6").Value = 181

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 105
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 346

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 266

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 158
$ws.Range("K6").Value = 162
$ws.Range("K7").Value = 466

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 389

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 66
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 90
$ws.Range("K3").Value = 92
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 333

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 53
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K6").Value = 294
$ws.Range("K7").Value = 1188

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 183
$ws.Range("K3").Value = 186
$ws.Range("K7").Value = 673

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 91
